# minor fixes and TODO added
# Appends newly tasted/researched rums to the bottom of the Rum Howler
# ratings sheet (Sheet1: A=Rum, B=Score, C=ReviewCount, D=Source).
# A handful of rows are TODO stubs - only the name is known yet, score /
# review-count / source will be filled in later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 381

$names = @(
    "Baltic Dark muovipullo",
    "Rhum Negrita Dark Signature",
    "Barra Oak Bay Dark",
    "Ron Cabana Blanco muovipullo",
    "Rhum Negrita White Signature",
    "Propeller Dark Rum",
    "Baltic Light muovipullo",
    "Rommiviina muovipullo",
    "Old Pascas Dark",
    "Barra Oak Bay Dark muovipullo",
    "Rhum Negrita Dark Signature muovipullo",
    "Barceló Blanco Añejado",
    "Ching Shih Dark Spiced",
    "Barbados Rommi muovipullo",
    "Barra Oak Bay White muovipullo",
    "Negrita Spiced Golden",
    "Ripa's Honey & Chili",
    "Planteray O.F.T.D",
    "Bacardi Carta Negra",
    "Stroh"
)

# Scores / review counts / sources, keyed by row number. Rows 393, 394 and
# 397 are TODO stubs - only the rum name is filled in for now.
$scores = [ordered]@{
    381 = 22;  382 = 29;  383 = 48;  384 = 54;  385 = 47;
    386 = 48;  387 = 20;  388 = 30;  389 = 49;  390 = 48;
    391 = 29;  392 = 52;  395 = 63;  396 = 47;  398 = 76;
    399 = 48;  400 = 39
}
$reviewCounts = [ordered]@{
    381 = 5;   382 = 95;  383 = 5;   384 = 16;  385 = 11;
    386 = 9;   387 = 1;   388 = 1;   389 = 31;  390 = 5;
    391 = 95;  392 = 25;  395 = 4;   396 = 24;  398 = 285;
    399 = 117; 400 = 32
}
$sources = [ordered]@{
    381 = "rumratings.com"; 382 = "rumratings.com"; 383 = "rumratings.com";
    384 = "isokaato.com";   385 = "rumratings.com"; 386 = "rumratings.com";
    387 = "rumratings.com"; 388 = "rumratings.com"; 389 = "rumratings.com";
    390 = "rumratings.com"; 391 = "rumratings.com"; 392 = "rumratings.com";
    395 = "isokaato.com";   396 = "rumratings.com";
    398 = "rumratings.com"; 399 = "rumratings.com"; 400 = "rumratings.com"
}

# Fill column A (Rum name) for every new row first, then B/C, then D -
# mirrors how the sheet was actually typed up (names jotted down first,
# ratings and sources filled in afterwards).
for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($firstNewRow + $i, 1).Value = $names[$i]
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 2).Value = $scores[$row]
}

foreach ($row in $reviewCounts.Keys) {
    $ws.Cells.Item($row, 3).Value = $reviewCounts[$row]
}

foreach ($row in $sources.Keys) {
    $ws.Cells.Item($row, 4).Value = $sources[$row]
}

# Leave the selection where editing left off, scrolled down to the new rows.
$ws.Range("D401").Select() | Out-Null
